$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 12plant_SNP
$ws2 = $wb.Worksheets.Item(2)   # Domest_SNP
$ws3 = $wb.Worksheets.Item(3)   # 12plant_gene
$ws4 = $wb.Worksheets.Item(4)   # Domest_gene

# ---------------------------------------------------------------------
# Sheet1 (12plant_SNP): convert D2:D13 sum formulas into one shared
# formula block, then add the new summary / "top SNPs" rows below.
# ---------------------------------------------------------------------
$ws1.Range("D2:D13").Formula = "=SUM(B2:C2)"

$ws1.Range("D14").Formula = "=SUM(D2:D13)"
$ws1.Range("E14").Formula = "=46000/D14"

# Row 17 + 18: blank, formatted cells (copy format from an already
# "plain Calibri" styled cell elsewhere in the workbook)
$ws3.Range("B7").Copy()
$ws1.Range("B17:G17").PasteSpecial(-4122)
$ws1.Range("B18:E18").PasteSpecial(-4122)
$ws1.Range("A19:B19").PasteSpecial(-4122)
$ws1.Range("D19:F19").PasteSpecial(-4122)

$ws1.Range("A18").Value = "Top 1000 SNPs only"

# Row 19: header-ish numbers, some styled like row 17/18, some not
$ws1.Range("A19").Value = 1
$ws1.Range("B19").Value = 2
$ws1.Range("C19").Value = 3
$ws1.Range("D19").Value = 4
$ws1.Range("E19").Value = 5
$ws1.Range("F19").Value = 6
$ws1.Range("G19").Value = 7
$ws1.Range("H19").Value = 8
$ws1.Range("I19").Value = 9
$ws1.Range("J19").Value = 10

# Row 20: raw data + row sum
$ws1.Range("A20").Value = 6564
$ws1.Range("B20").Value = 1264
$ws1.Range("C20").Value = 416
$ws1.Range("D20").Value = 205
$ws1.Range("E20").Value = 75
$ws1.Range("F20").Value = 24
$ws1.Range("G20").Value = 34
$ws1.Range("H20").Value = 8
$ws1.Range("I20").Value = 1
$ws1.Range("J20").Value = 1
$ws1.Range("K20").Formula = "=SUM(A20:J20)"

# Row 21: ratio
$ws1.Range("K21").Formula = "=A20/K20"

# Column C width: keep the same displayed width but mark it as an
# explicit (non "best fit") custom width, same as the author's edit.
$ws1.Columns.Item(3).ColumnWidth = 13.8

# Selection on sheet1 moved from A17 to E14
$ws1.Range("E14").Select()

# ---------------------------------------------------------------------
# Sheet3 (12plant_gene): add row total + ratio column N
# ---------------------------------------------------------------------
$ws3.Range("N7").Formula = "=SUM(B7:M7)"
$ws3.Range("N8").Formula = "=555/990"

$ws3.Range("G16").Select()

# ---------------------------------------------------------------------
# Active sheet moves from Domest_gene (4th tab) to 12plant_gene (3rd tab)
# ---------------------------------------------------------------------
$ws3.Activate()

Write-Output "done"
